$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Bands" sheet — tweak a few stat cells for the Veig's Band row.
# ---------------------------------------------------------------------------
$bands = $wb.Worksheets.Item("Bands")
$bands.Range("E2").Value = "VIGILANCE"     # faculty
$bands.Range("G2").Value = 18              # huntThreshold
$bands.Range("I2").Value = 10              # hopePts
$bands.Range("J2").Value = 1               # shadowPts

# ---------------------------------------------------------------------------
# 2. New "Heroes" sheet — inserted right after "Bands".
# ---------------------------------------------------------------------------
$heroes = $wb.Worksheets.Add($null, $bands)
$heroes.Name = "Heroes"

$heroes.Cells.Item(1, 1).Value = "id"
$heroes.Cells.Item(1, 2).Value = "name"
$heroes.Cells.Item(2, 1).Value = "qwer1234"
$heroes.Cells.Item(2, 2).Value = "Veig"

# ---------------------------------------------------------------------------
# 3. "Allies" sheet — content is unchanged, only the selected cell moves.
# ---------------------------------------------------------------------------
$allies = $wb.Worksheets.Item("Allies")

# ---------------------------------------------------------------------------
# 4. New "Items" sheet — inserted right after "Allies".
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Add($null, $allies)
$items.Name = "Items"

$items.Cells.Item(1, 1).Value = "id"
$items.Cells.Item(1, 2).Value = "idOwner"
$items.Cells.Item(1, 3).Value = "name"
$items.Cells.Item(1, 4).Value = "description"
$items.Cells.Item(1, 5).Value = "slot"
$items.Cells.Item(1, 6).Value = "benefit1"
$items.Cells.Item(1, 7).Value = "benefit2"
$items.Cells.Item(1, 8).Value = "isCursed"
$items.Cells.Item(1, 9).Value = "type"

$items.Cells.Item(2, 1).Value = "lol"
$items.Cells.Item(2, 2).Value = "qwer1234"
$items.Cells.Item(2, 3).Value = "Sword of Destiny"
$items.Cells.Item(2, 4).Value = "A legendary sword with a glowing blade."
$items.Cells.Item(2, 5).Value = "ItemSlotTypeTOR.WEAPON"
$items.Cells.Item(2, 6).Value = "SkillTypeTOR.BATTLE "
$items.Cells.Item(2, 7).Value = "SkillTypeTOR.NONE"
$items.Cells.Item(2, 8).Value = "'False"
$items.Cells.Item(2, 9).Value = "MagicItemType.UNUSUAL"

# ---------------------------------------------------------------------------
# 5. Selections for every sheet, matching the final saved view state.
# ---------------------------------------------------------------------------
$bands.Range("G8").Select()
$heroes.Range("A2").Select()
$allies.Range("M1").Select()
$items.Range("E12").Select()

# "Items" ends up the active tab.
$items.Activate()
